$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Val)
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.Style = "Normal"
}

Set-TextValue "D2" '59.342.04'
Set-TextValue "E2" '  +1.01%  '
Set-TextValue "D3" '2.606.20'
Set-TextValue "E3" '  +0.84%  '
Set-TextValue "E4" '  +0.01%  '
Set-TextValue "D5" '541.71'
Set-TextValue "E5" '  +4.21%  '
Set-TextValue "D6" '141.69'
Set-TextValue "E6" '  +1.78%  '
Set-TextValue "E7" '  +0.11%  '
Set-TextValue "E8" '  +0.26%  '
Set-TextValue "E10" '  +2.38%  '
Set-TextValue "E11" '  +2.02%  '
Set-TextValue "E12" '  +0.56%  '
Set-TextValue "D13" '3.060.04'
Set-TextValue "D14" '59.281.21'
Set-TextValue "E14" '  +0.96%  '
Set-TextValue "E15" '  +1.27%  '
Set-TextValue "E16" '  +1.15%  '
Set-TextValue "D17" '2.603.25'
Set-TextValue "E17" '  +1.46%  '
Set-TextValue "D18" '341.22'
Set-TextValue "E18" '  +0.85%  '
Set-TextValue "E19" '  +1.67%  '
Set-TextValue "D20" '10.15'
Set-TextValue "E20" '  +0.11%  '
Set-TextValue "E21" '  -1.39%  '
Set-TextValue "D22" '0.999'
Set-TextValue "E22" '  +0.06%  '
Set-TextValue "D23" '67.52'
Set-TextValue "E23" '  +1.85%  '
Set-TextValue "E24" '  +1.55%  '
Set-TextValue "E25" '  -1.43%  '
Set-TextValue "E26" '  +0.09%  '
Set-TextValue "D27" '7.25'
Set-TextValue "E27" '  +3.20%  '
Set-TextValue "D28" '0.0₃0747'
Set-TextValue "E28" '  +4.20%  '
Set-TextValue "D30" '1.68'
Set-TextValue "E30" '  +7.06%  '
Set-TextValue "E31" '  -1.91%  '
Set-TextValue "D32" '18.76'
Set-TextValue "E32" '  +0.08%  '
Set-TextValue "D33" '149.84'
Set-TextValue "E33" '  +0.88%  '
Set-TextValue "E34" '  +0.54%  '
Set-TextValue "E35" '  +0.10%  '
Set-TextValue "D36" '37.22'
Set-TextValue "E36" '  +2.02%  '
Set-TextValue "E37" '  +0.77%  '
Set-TextValue "D38" '0.837'
Set-TextValue "E38" '  +1.29%  '
Set-TextValue "E39" '  +1.70%  '
Set-TextValue "E40" '  +1.88%  '
Set-TextValue "E41" '  +0.12%  '
Set-TextValue "D42" '275.94'
Set-TextValue "E42" '  +0.45%  '
Set-TextValue "D43" '0.601'
Set-TextValue "E43" '  +1.93%  '
Set-TextValue "D44" '10.72'
Set-TextValue "E44" '  -0.21%  '
Set-TextValue "D45" '0.0955'
Set-TextValue "E45" '  +0.79%  '
Set-TextValue "D46" '0.0526'
Set-TextValue "E46" '  +1.16%  '
Set-TextValue "D47" '1.954.52'
Set-TextValue "E47" '  -1.19%  '
Set-TextValue "D48" '18.57'
Set-TextValue "E48" '  +3.94%  '
Set-TextValue "E49" '  +1.93%  '
Set-TextValue "E50" '  +0.46%  '
Set-TextValue "D51" '110.93'
Set-TextValue "E51" '  -0.90%  '
